$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wealth Class in Allocation Row: 12 -> 22
$ws.Range("B2").Value = 22

# Wealth Class in Cash Flow Row: 16 -> 20
$ws.Range("B3").Value = 20

# Wealth Row: 7 -> 9
$ws.Range("B5").Value = 9

# Move the active selection to B3 (test for reduction of wealth class)
$ws.Range("B3").Select()
